$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-01-11 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-01-12 Thursday", 2) | Out-Null

# Update each math-expression cell in the table by position (row, col)
# to avoid any ambiguity from repeated/substring text values.
$t = $d.Tables(1)

$newValues = @(
    @("83-22=", "54-25=", "35+12=", "99-93=", "26+6="),
    @("43-41=", "43-27=", "26+56=", "39+4=", "83-60="),
    @("66-30=", "7+21=", "70+27=", "88-1=", "77-56="),
    @("58-48=", "46+45=", "65-46=", "54-34=", "28+47="),
    @("46-12=", "6+41=", "91-89=", "81-52=", "65+26="),
    @("41+34=", "81-70=", "42+33=", "11+54=", "9+47="),
    @("53+41=", "63-55=", "86-75=", "59+28=", "66-17="),
    @("15-13=", "92-5=", "94-73=", "85-53=", "20+65="),
    @("51+18=", "99-82=", "58+17=", "74-55=", "89+4="),
    @("76-46=", "7+83=", "24-2=", "62-53=", "73-66="),
    @("76-20=", "16-11=", "99-93=", "17+45=", "60+2="),
    @("98-43=", "80-7=", "44-10=", "1+54=", "11+70="),
    @("18-12=", "6+42=", "99-85=", "17-5=", "5+73="),
    @("56-8=", "88-79=", "31-26=", "79-33=", "53-49="),
    @("12+53=", "14+52=", "12+69=", "4-4=", "64-29="),
    @("61-48=", "54+38=", "86-58=", "4+24=", "35-4="),
    @("57-6=", "26+20=", "28+66=", "9+23=", "25+42="),
    @("47-3=", "1+58=", "84+14=", "49+15=", "26-12="),
    @("1+5=", "23-20=", "78-72=", "97-93=", "36+0="),
    @("10+28=", "92-68=", "68+17=", "83-38=", "14-11=")
)

for ($r = 1; $r -le $newValues.Length; $r++) {
    $rowValues = $newValues[$r - 1]
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $rowValues[$c - 1]
    }
}

Write-Output "done"